# Revision 17 Ago USA
# Adds the "STATUS BY PO" / "STATUS BY JOB" reference block to Sheet2
# and refreshes the selections on Sheet1 / Sheet2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- New reference data block on Sheet2 (rows 13-23) -----------------
# Entry order below reproduces the shared-strings insertion order of the
# authored workbook.
$ws2.Range("C13").Value = "STATUS BY PO"
$ws2.Range("C14").Value = "OPEN"
$ws2.Range("C16").Value = "COMPLETED"
$ws2.Range("C19").Value = "STATUS BY JOB"
$ws2.Range("C15").Value = "HOLD"
$ws2.Range("D13").Value = "COMMENTS"
$ws2.Range("C21").Value = "ACTIVE"
$ws2.Range("D22").Value = "All POS are woks status"
$ws2.Range("C17").Value = "CLOSED"
$ws2.Range("C20").Value = "OPEN"
$ws2.Range("C22").Value = "HOLD"
$ws2.Range("C23").Value = "CLOSED"

# --- Selection refresh -------------------------------------------------
# Sheet1's cursor moved to A7 ...
$ws1.Range("A7").Select()
# ... while Sheet2 stays the active/visible tab, cursor parked on C20.
$ws2.Activate()
$ws2.Range("C20").Select()
